# Update cryptos price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "25.950.62"
$ws.Range("E2").Value = "  -0.45%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.621.38"
$ws.Range("E3").Value = "  -1.12%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.48%  "

# Row 5: BNB
$ws.Range("D5").Value = "'212.84"
$ws.Range("E5").Value = "  -0.86%  "

# Row 6: XRP
$ws.Range("D6").Value = "'0.498"
$ws.Range("E6").Value = "  -1.31%  "

# Row 7: USDC
$ws.Range("E7").Value = "  -0.46%  "

# Row 8: Cardano
$ws.Range("E8").Value = "  -0.70%  "

# Row 9: Dogecoin
$ws.Range("E9").Value = "  -1.22%  "

# Row 10: Solana
$ws.Range("D10").Value = "'18.44"
$ws.Range("E10").Value = "  -1.33%  "

# Row 11: TRON
$ws.Range("E11").Value = "  -0.30%  "

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.845.86"
$ws.Range("E12").Value = "  -1.15%  "

# Row 13: WrappedEther
$ws.Range("D13").Value = "1.612.44"
$ws.Range("E13").Value = "  -1.92%  "

# Row 14: Polkadot
$ws.Range("D14").Value = "'4.15"
$ws.Range("E14").Value = "  -1.69%  "

# Row 15: Polygon
$ws.Range("D15").Value = "'0.526"
$ws.Range("E15").Value = "  -1.21%  "

# Row 16: WrappedBTC
$ws.Range("D16").Value = "25.958.46"
$ws.Range("E16").Value = "  -0.43%  "

# Row 17: Litecoin
$ws.Range("D17").Value = "'61.72"
$ws.Range("E17").Value = "  -1.11%  "

# Row 18: ShibaInu
$ws.Range("D18").Value = "0.0₃0739"
$ws.Range("E18").Value = "  -1.42%  "

# Row 19: Dai
$ws.Range("E19").Value = "  -0.52%  "

# Row 20: BitcoinCash
$ws.Range("D20").Value = "'192.57"
$ws.Range("E20").Value = "  +0.77%  "

# Row 21: Uniswap
$ws.Range("D21").Value = "'4.26"
$ws.Range("E21").Value = "  -0.48%  "

# Row 22: Avalanche
$ws.Range("D22").Value = "'9.54"
$ws.Range("E22").Value = "  -0.57%  "

# Row 23: Chainlink
$ws.Range("E23").Value = "  -2.02%  "

# Row 24: Stellar
$ws.Range("E24").Value = "  +1.92%  "

# Row 25: Monero
$ws.Range("D25").Value = "'143.80"

# Row 26: BinanceUSD
$ws.Range("E26").Value = "  -0.44%  "

# Row 27: Toncoin
$ws.Range("D27").Value = "'1.72"
$ws.Range("E27").Value = "  -2.67%  "

# Row 28: Cosmos
$ws.Range("D28").Value = "'6.65"

# Row 29: EthereumClassic
$ws.Range("D29").Value = "'15.26"
$ws.Range("E29").Value = "  +0.13%  "

# Row 30: PancakeSwap
$ws.Range("E30").Value = "  -0.89%  "

# Row 31: Hedera
$ws.Range("D31").Value = "'0.0480"
$ws.Range("E31").Value = "  -1.38%  "

# Row 32: Filecoin
$ws.Range("D32").Value = "'3.13"
$ws.Range("E32").Value = "  -1.55%  "

# Row 33: InternetComputer(DFINITY)
$ws.Range("E33").Value = "  -2.50%  "

# Row 34: LidoDAOToken
$ws.Range("D34").Value = "'1.50"
$ws.Range("E34").Value = "  -0.90%  "

# Row 35: HuobiToken
$ws.Range("E35").Value = "  -1.33%  "

# Row 36: Maker
$ws.Range("D36").Value = "1.127.20"
$ws.Range("E36").Value = "  +0.24%  "

# Row 37: ARBITRUM
$ws.Range("D37").Value = "'0.845"
$ws.Range("E37").Value = "  -3.86%  "

# Row 38: MXToken
$ws.Range("D38").Value = "'2.40"
$ws.Range("E38").Value = "  -2.47%  "

# Row 39: ImmutableX
$ws.Range("D39").Value = "'0.516"
$ws.Range("E39").Value = "  -1.60%  "

# Row 40: VeChain
$ws.Range("E40").Value = "  -1.03%  "

# Row 41: Quant
$ws.Range("D41").Value = "'97.72"
$ws.Range("E41").Value = "  -1.15%  "

# Row 42: RocketPoolETH
$ws.Range("D42").Value = "1.757.09"
$ws.Range("E42").Value = "  -0.93%  "

# Row 43: TrustWalletToken
$ws.Range("D43").Value = "'0.756"
$ws.Range("E43").Value = "  -3.86%  "

# Row 44: FraxShare
$ws.Range("D44").Value = "'5.09"
$ws.Range("E44").Value = "  -4.00%  "

# Row 45: BabyDogeCoin
$ws.Range("E45").Value = "  -0.55%  "

# Row 46: RenderToken
$ws.Range("E46").Value = "  +2.22%  "

# Row 47: Aave
$ws.Range("D47").Value = "'54.18"
$ws.Range("E47").Value = "  -1.99%  "

# Row 48: Cronos
$ws.Range("E48").Value = "  -1.90%  "

# Row 49: Mantle
$ws.Range("D49").Value = "'0.410"
$ws.Range("E49").Value = "  -1.07%  "

# Row 50: EnergySwap
$ws.Range("D50").Value = "'7.48"
$ws.Range("E50").Value = "  -1.36%  "

# Row 51: USDD
$ws.Range("E51").Value = "  -0.24%  "
